# Validation set.xlsx - update
# 1. Rename "Original dataset" sheet -> "original dataset" (defined names follow automatically)
# 2. Fix the two orphaned #REF! defined names (no3_, po4_) so they keep the sheet-qualified prefix
# 3. Move the active/selected tab from "dataset used for testing ML" to "original dataset"
# 4. Update the saved selections on each sheet (F36 on original dataset, AF14 stays on the ML sheet)

$wb = $excel.ActiveWorkbook

# --- Rename the first sheet -------------------------------------------------
$ws1 = $wb.Worksheets.Item("Original dataset")
$ws1.Name = "original dataset"

# Re-point the two #REF! defined names so they keep the sheet-qualified prefix
# (Excel drops it on a plain sheet rename; restore it to match the sheet rename semantics).
$wb.Names.Item("no3_").RefersTo = "='original dataset'!#REF!"
$wb.Names.Item("po4_").RefersTo = "='original dataset'!#REF!"

# --- Move the active tab / selection ----------------------------------------
$ws2 = $wb.Worksheets.Item("dataset used for testing ML")

# Make "original dataset" the active sheet/tab and select F36 on it.
$ws1.Activate()
$ws1.Range("F36").Select()

# Keep the existing selection on the other sheet (AF14), just no longer the active tab.
$ws2.Range("AF14").Select()

# Restore the originally active sheet as the active one again so the saved
# file's tabSelected/active tab reflect "original dataset".
$ws1.Activate()
